$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "treatments" sheet: the "Fenced" treatment row (row 5) had its first/current
# application dates revised from "NA" to "TBD".
# ---------------------------------------------------------------------------
$wsTreat = $wb.Worksheets.Item("treatments")
$wsTreat.Cells.Item(5,4).Value = "TBD"   # D5 first_appl
$wsTreat.Cells.Item(5,5).Value = "TBD"   # E5 curr_appl

# ---------------------------------------------------------------------------
# "plan" sheet: plot-type determination revised to account for fenced plots.
# This reshuffles several n/p/k/exclose (columns E-H) values and fills in a
# new "note_plan" column (K) flagging plots as "not core plot" or
# "core control" where applicable.
# ---------------------------------------------------------------------------
$wsPlan = $wb.Worksheets.Item("plan")

$wsPlan.Cells.Item(2,11).Value = "not core plot"

$wsPlan.Cells.Item(3,5).Value = "Nitrogen"
$wsPlan.Cells.Item(3,6).Value = "Phosphorus"
$wsPlan.Cells.Item(3,7).Value = "Potassium"
$wsPlan.Cells.Item(3,8).Value = "Fence"

$wsPlan.Cells.Item(4,5).Value = "Nitrogen"

$wsPlan.Cells.Item(5,5).Value = "Control"
$wsPlan.Cells.Item(5,6).Value = "Control"

$wsPlan.Cells.Item(6,5).Value = "Nitrogen"

$wsPlan.Cells.Item(7,11).Value = "not core plot"

$wsPlan.Cells.Item(8,11).Value = "not core plot"

$wsPlan.Cells.Item(9,5).Value = "Control"
$wsPlan.Cells.Item(9,7).Value = "Control"
$wsPlan.Cells.Item(9,11).Value = "core control"

$wsPlan.Cells.Item(10,6).Value = "Control"

$wsPlan.Cells.Item(11,11).Value = "not core plot"

$wsPlan.Cells.Item(12,5).Value = "Nitrogen"
$wsPlan.Cells.Item(12,6).Value = "Phosphorus"

$wsPlan.Cells.Item(13,5).Value = "Control"
$wsPlan.Cells.Item(13,8).Value = "Fence"

$wsPlan.Cells.Item(14,6).Value = "Phosphorus"
$wsPlan.Cells.Item(14,7).Value = "Potassium"

$wsPlan.Cells.Item(16,8).Value = "Fence"

$wsPlan.Cells.Item(17,5).Value = "Nitrogen"
$wsPlan.Cells.Item(17,6).Value = "Phosphorus"

$wsPlan.Cells.Item(18,7).Value = "Potassium"

$wsPlan.Cells.Item(19,5).Value = "Nitrogen"
$wsPlan.Cells.Item(19,6).Value = "Control"

$wsPlan.Cells.Item(20,11).Value = "not core plot"

$wsPlan.Cells.Item(21,7).Value = "Potassium"

$wsPlan.Cells.Item(22,5).Value = "Control"

$wsPlan.Cells.Item(23,11).Value = "not core plot"

$wsPlan.Cells.Item(24,6).Value = "Phosphorus"
$wsPlan.Cells.Item(24,7).Value = "Control"

$wsPlan.Cells.Item(25,7).Value = "Potassium"
$wsPlan.Cells.Item(25,8).Value = "Fence"

$wsPlan.Cells.Item(26,11).Value = "not core plot"

$wsPlan.Cells.Item(27,11).Value = "not core plot"

$wsPlan.Cells.Item(28,5).Value = "Nitrogen"

$wsPlan.Cells.Item(29,5).Value = "Control"
$wsPlan.Cells.Item(29,7).Value = "Control"
$wsPlan.Cells.Item(29,11).Value = "core control"

$wsPlan.Cells.Item(30,6).Value = "Phosphorus"

$wsPlan.Cells.Item(31,5).Value = "Control"
$wsPlan.Cells.Item(31,6).Value = "Phosphorus"

$wsPlan.Cells.Item(32,11).Value = "not core plot"

$wsPlan.Cells.Item(33,5).Value = "Nitrogen"

$wsPlan.Cells.Item(34,5).Value = "Nitrogen"
$wsPlan.Cells.Item(34,6).Value = "Phosphorus"
$wsPlan.Cells.Item(34,7).Value = "Control"

$wsPlan.Cells.Item(35,5).Value = "Nitrogen"
$wsPlan.Cells.Item(35,8).Value = "Fence"

$wsPlan.Cells.Item(36,6).Value = "Control"
$wsPlan.Cells.Item(36,8).Value = "Fence"

$wsPlan.Cells.Item(37,5).Value = "Control"
$wsPlan.Cells.Item(37,6).Value = "Control"
$wsPlan.Cells.Item(37,7).Value = "Potassium"

$wsPlan.Cells.Item(38,5).Value = "Control"
$wsPlan.Cells.Item(38,6).Value = "Control"
$wsPlan.Cells.Item(38,7).Value = "Control"
$wsPlan.Cells.Item(38,11).Value = "core control"

$wsPlan.Cells.Item(39,11).Value = "not core plot"

$wsPlan.Cells.Item(40,6).Value = "Phosphorus"
$wsPlan.Cells.Item(40,7).Value = "Potassium"

$wsPlan.Cells.Item(41,11).Value = "not core plot"

$wsPlan.Cells.Item(42,11).Value = "not core plot"

$wsPlan.Cells.Item(43,5).Value = "Nitrogen"
$wsPlan.Cells.Item(43,7).Value = "Potassium"

# ---------------------------------------------------------------------------
# Window / selection state: "plan" sheet selection moves to G43, and the
# "treatments" sheet becomes the active tab with selection on E6 (the "block"
# sheet, previously active, loses its tabSelected flag as a side effect).
# ---------------------------------------------------------------------------
$wsPlan.Activate()
$wsPlan.Range("G43").Select()

$wsTreat.Activate()
$wsTreat.Range("E6").Select()
